$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 130; existing rows 130-134 shift down to 131-135
$ws.Rows.Item(130).Insert()

# Populate the newly inserted row 130 with the new weekly price record
$ws.Cells.Item(130, 1).Value = 10
$ws.Cells.Item(130, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(130, 3).Value = "La Araucanía"
$ws.Cells.Item(130, 4).Value = 45267
$ws.Cells.Item(130, 5).Value = 9
$ws.Cells.Item(130, 6).Value = 100112022
$ws.Cells.Item(130, 7).Value = "Arveja Verde"
$ws.Cells.Item(130, 8).Value = "Sin especificar"
$ws.Cells.Item(130, 9).Value = "Primera"
$ws.Cells.Item(130, 10).Value = 80
$ws.Cells.Item(130, 11).Value = 25000
$ws.Cells.Item(130, 12).Value = 25000
$ws.Cells.Item(130, 13).Value = 25000
$ws.Cells.Item(130, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(130, 15).Value = "Región del Maule"
$ws.Cells.Item(130, 16).Value = 1000
$ws.Cells.Item(130, 17).Value = 25
$ws.Cells.Item(130, 18).Value = "Hortaliza"

# Ensure date number format (s="2", the datetime number format) carries over to the new D130 cell
$ws.Cells.Item(130, 4).NumberFormat = $ws.Cells.Item(131, 4).NumberFormat
